$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new line records (line7, line8) are inserted right after the existing
# line6 row, pushing the 8 extr* rows down by two. This is a plain row
# insert at row 8 followed by filling in the new rows and touching up the
# few numeric values that changed further down.
$ws.Rows("8:9").Insert()

# Match the header-column style (bold font, thin border, centered) that the
# inserted rows should carry in column A (same as all the other data rows).
$ws.Range("A8:A9").Font.Bold = $true
$ws.Range("A8:A9").HorizontalAlignment = -4108
$ws.Range("A8:A9").VerticalAlignment = -4160
$ws.Range("A8:A9").Borders.LineStyle = 1
$ws.Range("A8:A9").Borders.Weight = 2

# --- New row 8: line7 ---
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# --- New row 9: line8 ---
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# --- Row 10 (was extr1): only the running index changes ---
$ws.Range("A10").Value = 8

# --- Row 11 (was extr2): index + in_service flips to TRUE ---
$ws.Range("A11").Value = 9
$ws.Range("E11").Value = $true

# --- Row 12 (was extr3): index + in_service flips to TRUE ---
$ws.Range("A12").Value = 10
$ws.Range("E12").Value = $true

# --- Row 13 (was extr4): only the running index changes ---
$ws.Range("A13").Value = 11

# --- Row 14 (was extr5): only the running index changes ---
$ws.Range("A14").Value = 12

# --- Row 15 (was extr6): index + in_service flips to FALSE ---
$ws.Range("A15").Value = 13
$ws.Range("E15").Value = $false

# --- Row 16 (was extr7): index + in_service flips to TRUE ---
$ws.Range("A16").Value = 14
$ws.Range("E16").Value = $true

# --- Row 17 (was extr8): only the running index changes ---
$ws.Range("A17").Value = 15
